# This script removes the "Empirical multiplicity of solid diffusivity" (mD)
# parameter row from the "Parameters" sheet, fixing the linEIS regression
# referenced in the commit message. Deleting the row shifts every
# subsequent row up by one (old row 64 -> new row 63, ..., old row 75 ->
# new row 74) and causes the three now-unused shared strings
# ("Empirical multiplicity of solid diffusivity", "mD", "m_\mathrm{D}")
# to be dropped automatically when the workbook is saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")

# Row 63 currently holds the mD parameter definition:
#   B63 = "mD", C63 = "Empirical multiplicity of solid diffusivity",
#   D63 = "m_\mathrm{D}", E63 = 1, F63 = 0, G63 = "s"
# Deleting the entire row removes it and shifts rows 64:75 up to 63:74.
$ws.Rows.Item(63).Delete()
